$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'Pipeline(steps=[(''scaler'', RobustScaler()),
                (''selector'',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'',
                                                     random_state=42))),
                (''model'',
                 AdaBoostClassifier(estimator=LogisticRegression(C=5,
                                                                 max_iter=1000,
                                                                 penalty=''l1'',
                                                                 random_state=42,
                                                                 solver=''saga''),
                                    n_estimators=10, random_state=42))])'
$ws.Range("B2").Value = 0.7499999999999999
$ws.Range("C2").Value = '{''selector'': SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'', random_state=42)), ''scaler'': RobustScaler(), ''model__n_estimators'': 10, ''model__estimator__solver'': ''saga'', ''model__estimator__penalty'': ''l1'', ''model__estimator__class_weight'': None, ''model__estimator__C'': 5}'
$ws.Range("D2").Value = 0.5723676963620047
$ws.Range("E2").Value = 0.5215124262139187
$ws.Range("F2").Value = 0.8
$ws.Range("G2").Value = 0.5334805720736272
$ws.Range("H2").Value = 0.4847654584221748
$ws.Range("I2").Value = 0.6666666666666666
$ws.Range("J2").Value = 0.6596379803112099
$ws.Range("K2").Value = 0.6032338308457712
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = '[1 0 1 1 1 1 0 1 0 1 0 1 0 1 1 0 0 1 1 1 1 0 1 1]'
$ws.Range("N2").Value = '[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]'

# Row 3
$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', StandardScaler()), (''selector'', None),
                (''model'',
                 AdaBoostClassifier(estimator=LogisticRegression(C=0.0001,
                                                                 max_iter=1000,
                                                                 penalty=''l1'',
                                                                 random_state=42,
                                                                 solver=''saga''),
                                    n_estimators=5, random_state=42))])'
$ws.Range("B3").Value = 0.7499999999999999
$ws.Range("C3").Value = '{''selector'': None, ''scaler'': StandardScaler(), ''model__n_estimators'': 5, ''model__estimator__solver'': ''saga'', ''model__estimator__penalty'': ''l1'', ''model__estimator__class_weight'': None, ''model__estimator__C'': 0.0001}'
$ws.Range("D3").Value = 0.4773453735546231
$ws.Range("E3").Value = 0.4405008495777726
$ws.Range("F3").Value = 0.8
$ws.Range("G3").Value = 0.4315136760761198
$ws.Range("H3").Value = 0.4248754578754579
$ws.Range("I3").Value = 0.6666666666666666
$ws.Range("J3").Value = 0.5598690671031098
$ws.Range("K3").Value = 0.4994871794871795
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = '[1 1 0 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0]'
$ws.Range("N3").Value = '[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]'

# Row 4
$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7fa034091910>),
                (''model'',
                 AdaBoostClassifier(estimator=LogisticRegression(C=0.001,
                                                                 class_weight=''balanced'',
                                                                 max_iter=1000,
                                                                 random_state=42,
                                                                 solver=''liblinear''),
                                    n_estimators=5, random_state=42))])'
$ws.Range("B4").Value = 0.6959706959706959
$ws.Range("C4").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7fa03404f640>, ''scaler'': MinMaxScaler(), ''model__n_estimators'': 5, ''model__estimator__solver'': ''liblinear'', ''model__estimator__penalty'': ''l2'', ''model__estimator__class_weight'': ''balanced'', ''model__estimator__C'': 0.001}'
$ws.Range("D4").Value = 0.5474990330408427
$ws.Range("E4").Value = 0.4933154927509767
$ws.Range("F4").Value = 0.8837209302325582
$ws.Range("G4").Value = 0.5108881358712153
$ws.Range("H4").Value = 0.4529678699436764
$ws.Range("I4").Value = 0.7916666666666666
$ws.Range("J4").Value = 0.6254838709677418
$ws.Range("K4").Value = 0.5909677419354837
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = '[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1]'
$ws.Range("N4").Value = '[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]'

# Row 5
$ws.Range("A5").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'',
                                                     random_state=42))),
                (''model'',
                 AdaBoostClassifier(estimator=LogisticRegression(C=3,
                                                                 max_iter=1000,
                                                                 random_state=42,
                                                                 solver=''saga''),
                                    n_estimators=10, random_state=42))])'
$ws.Range("B5").Value = 0.7669230769230768
$ws.Range("C5").Value = '{''selector'': SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'', random_state=42)), ''scaler'': MinMaxScaler(), ''model__n_estimators'': 10, ''model__estimator__solver'': ''saga'', ''model__estimator__penalty'': ''l2'', ''model__estimator__class_weight'': None, ''model__estimator__C'': 3}'
$ws.Range("D5").Value = 0.5072585339700348
$ws.Range("E5").Value = 0.4697612451960277
$ws.Range("F5").Value = 0.7368421052631579
$ws.Range("G5").Value = 0.4731771707371918
$ws.Range("H5").Value = 0.4285127674258109
$ws.Range("I5").Value = 0.5833333333333334
$ws.Range("J5").Value = 0.5838805087252292
$ws.Range("K5").Value = 0.555072463768116
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = '[0 1 1 0 0 1 0 0 0 0 1 1 1 0 0 1 1 0 1 1 1 1 1 1]'
$ws.Range("N5").Value = '[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]'

# Row 6
$ws.Range("A6").Value = 'Pipeline(steps=[(''scaler'', RobustScaler()),
                (''selector'',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                (''model'',
                 AdaBoostClassifier(estimator=LogisticRegression(C=5,
                                                                 max_iter=1000,
                                                                 penalty=''l1'',
                                                                 random_state=42,
                                                                 solver=''saga''),
                                    n_estimators=10, random_state=42))])'
$ws.Range("B6").Value = 0.7499999999999999
$ws.Range("C6").Value = '{''selector'': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), ''scaler'': RobustScaler(), ''model__n_estimators'': 10, ''model__estimator__solver'': ''saga'', ''model__estimator__penalty'': ''l1'', ''model__estimator__class_weight'': None, ''model__estimator__C'': 5}'
$ws.Range("D6").Value = 0.5642639621471427
$ws.Range("E6").Value = 0.4939514628497679
$ws.Range("F6").Value = 0.6285714285714286
$ws.Range("G6").Value = 0.5179437738360022
$ws.Range("H6").Value = 0.4465812483185364
$ws.Range("I6").Value = 0.4583333333333333
$ws.Range("J6").Value = 0.640677966101695
$ws.Range("K6").Value = 0.5830508474576273
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = '[1 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 0 0 1 1]'
$ws.Range("N6").Value = '[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]'
